$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1:C1").ColumnWidth = 36.28515625
